$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the refreshed cryptocurrency market snapshot.
# Columns B/C/D are forced to Text format before assignment so that
# numeric-looking strings (e.g. "216.01", "27.241.56") are preserved
# exactly as text instead of being coerced into floating point numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.241.56'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.688.47'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.01'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.01'
$ws.Range('E8').Value = '  +13.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.262'
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0627'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.927.62'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.686.62'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.19'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.32'
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.251.16'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.58'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.18'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0747'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.64'
$ws.Range('E23').Value = '  +4.73%  '
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.24'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.31'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.56'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0502'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.581.45'
$ws.Range('E32').Value = '  +6.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.41'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.26'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.70'
$ws.Range('E35').Value = '  +0.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.955'
$ws.Range('E36').Value = '  +5.96%  '
$ws.Range('E37').Value = '  +3.37%  '
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.07'
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.95'
$ws.Range('E41').Value = '  +3.67%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.64'
$ws.Range('E43').Value = '  -3.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.27'
$ws.Range('E44').Value = '  -2.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.835.62'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.786'
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '91.27'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('E48').Value = '  +5.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0106'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.17'
$ws.Range('E51').Value = '  +6.05%  '
